$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row: add D1 "CanAddSubHead" ---
$ws.Range("D1").Value = "CanAddSubHead"

# Give D1 the same border/format treatment as C1 (border-applied style),
# matching the look of the rest of the header row.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4163)  # xlPasteFormats

# --- Column C17 text changes from old description-style text to the new PropertyName "SelfAssessmentTax" ---
$ws.Range("C17").Value = "SelfAssessmentTax"

# --- New column D (CanAddSubHead) values, row 2-24 ---
$canAddSubHead = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 1
    6  = 1
    7  = 1
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 0
    22 = 0
    23 = 0
    24 = 0
}

foreach ($r in 2..24) {
    $ws.Range("D$r").Value = $canAddSubHead[$r]
}

# Apply the same format as column C (border-applied style) to the new D2:D24 values
$ws.Range("C2").Copy()
$ws.Range("D2:D24").PasteSpecial(-4163)  # xlPasteFormats

# --- New column E: the INSERT statement formula, now including PropertyName & CanAddSubHead ---
$formula = '=CONCATENATE("INSERT INTO [dbo].[ITHeadMaster] ([ExcelSrNo],[Description],[PropertyName],[CanAddSubHead],[Active],[AddedBy],[AddedDate]) VALUES (''",A2,"'',''",B2,"'',''",C2,,"'',",D2,",1,1,GETDATE())")'
$ws.Range("E2").Formula = $formula

$sharedFormula = '=CONCATENATE("INSERT INTO [dbo].[ITHeadMaster] ([ExcelSrNo],[Description],[PropertyName],[CanAddSubHead],[Active],[AddedBy],[AddedDate]) VALUES (''",A3,"'',''",B3,"'',''",C3,,"'',",D3,",1,1,GETDATE())")'
$ws.Range("E3:E24").Formula = $sharedFormula

# --- Column widths: C and D should share the same width (21) ---
$ws.Columns("D").ColumnWidth = $ws.Columns("C").ColumnWidth

# --- Selection matches the new target state ---
$ws.Range("D8").Select()
